$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: values are prefixed with a leading apostrophe so Excel stores them
# as text (matching the original inline-string cells) instead of converting
# numeric-looking strings (e.g. "0.05350" or "26.193.36") into numbers, which
# would silently drop significant trailing/leading zeros and re-interpret the
# "." thousands separators used in the price column.

$ws.Range("D2").Value = "'26.193.36"
$ws.Range("E2").Value = "'  +1.22%  "
$ws.Range("D3").Value = "'1.753.80"
$ws.Range("E3").Value = "'  +0.75%  "
$ws.Range("D4").Value = "'0.9983"
$ws.Range("E4").Value = "'  -0.08%  "
$ws.Range("D5").Value = "'238.80"
$ws.Range("E5").Value = "'  +4.70%  "
$ws.Range("D6").Value = "'0.9993"
$ws.Range("E6").Value = "'  +0.07%  "
$ws.Range("D7").Value = "'0.5301"
$ws.Range("E7").Value = "'  +3.27%  "
$ws.Range("D8").Value = "'0.2828"
$ws.Range("E8").Value = "'  +0.61%  "
$ws.Range("D9").Value = "'0.06214"
$ws.Range("E9").Value = "'  +2.05%  "
$ws.Range("D10").Value = "'1.747.27"
$ws.Range("E10").Value = "'  +0.31%  "
$ws.Range("D11").Value = "'0.07212"
$ws.Range("E11").Value = "'  +3.28%  "
$ws.Range("D12").Value = "'15.63"
$ws.Range("E12").Value = "'  +2.16%  "
$ws.Range("D13").Value = "'0.6515"
$ws.Range("E13").Value = "'  +2.55%  "
$ws.Range("D14").Value = "'4.652"
$ws.Range("E14").Value = "'  +3.48%  "
$ws.Range("D15").Value = "'78.95"
$ws.Range("E15").Value = "'  +3.32%  "
$ws.Range("D16").Value = "'0.9989"
$ws.Range("E16").Value = "'  -0.01%  "
$ws.Range("D17").Value = "'0.9976"
$ws.Range("E17").Value = "'  -0.14%  "
$ws.Range("D18").Value = "'26.068.91"
$ws.Range("E18").Value = "'  +0.73%  "
$ws.Range("D19").Value = "'11.86"
$ws.Range("E19").Value = "'  +3.14%  "
$ws.Range("D20").Value = "'0.000006772"
$ws.Range("E20").Value = "'  +2.83%  "
$ws.Range("D21").Value = "'1.970.65"
$ws.Range("E21").Value = "'  +0.15%  "
$ws.Range("D22").Value = "'4.351"
$ws.Range("E22").Value = "'  +6.16%  "
$ws.Range("D23").Value = "'8.794"
$ws.Range("E23").Value = "'  +3.82%  "
$ws.Range("D24").Value = "'5.269"
$ws.Range("E24").Value = "'  +2.92%  "
$ws.Range("D25").Value = "'139.60"
$ws.Range("E25").Value = "'  +0.12%  "
$ws.Range("D26").Value = "'1.517"
$ws.Range("D27").Value = "'15.43"
$ws.Range("E27").Value = "'  +2.78%  "
$ws.Range("D28").Value = "'1.830"
$ws.Range("E28").Value = "'  +0.96%  "
$ws.Range("D29").Value = "'105.50"
$ws.Range("E29").Value = "'  +2.91%  "
$ws.Range("D30").Value = "'0.08332"
$ws.Range("E30").Value = "'  +0.28%  "
$ws.Range("D31").Value = "'3.839"
$ws.Range("E31").Value = "'  +6.11%  "
$ws.Range("D32").Value = "'3.677"
$ws.Range("E32").Value = "'  +7.96%  "
$ws.Range("D33").Value = "'0.04635"
$ws.Range("E33").Value = "'  +5.73%  "
$ws.Range("D34").Value = "'2.652"
$ws.Range("E34").Value = "'  +1.28%  "
$ws.Range("D35").Value = "'1.029"
$ws.Range("E35").Value = "'  +6.08%  "
$ws.Range("D36").Value = "'0.6401"
$ws.Range("E36").Value = "'  +5.56%  "
$ws.Range("E37").Value = "'  +1.48%  "
$ws.Range("D38").Value = "'0.01625"
$ws.Range("E38").Value = "'  +4.32%  "
$ws.Range("D39").Value = "'2.003"
$ws.Range("E39").Value = "'  +4.39%  "
$ws.Range("D40").Value = "'0.9981"
$ws.Range("E40").Value = "'  +0.03%  "
$ws.Range("D41").Value = "'102.48"
$ws.Range("E41").Value = "'  +1.85%  "
$ws.Range("D42").Value = "'0.3980"
$ws.Range("E42").Value = "'  +3.78%  "
$ws.Range("D43").Value = "'0.7542"
$ws.Range("E43").Value = "'  +4.34%  "
$ws.Range("D44").Value = "'5.061"
$ws.Range("E44").Value = "'  +2.60%  "
$ws.Range("D45").Value = "'0.1160"
$ws.Range("E45").Value = "'  +4.85%  "
$ws.Range("D46").Value = "'6.447"
$ws.Range("E46").Value = "'  +2.00%  "
$ws.Range("D47").Value = "'0.05350"
$ws.Range("E47").Value = "'  -1.73%  "

# Row 48 and 49 swap: Aave <-> Elrond
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "'31.16"
$ws.Range("E48").Value = "'  +4.80%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'54.67"
$ws.Range("E49").Value = "'  +4.11%  "

$ws.Range("D50").Value = "'0.3510"
$ws.Range("E50").Value = "'  +3.67%  "
$ws.Range("D51").Value = "'7.624"
$ws.Range("E51").Value = "'  +1.37%  "
